# Weekly price-sheet update: a new weekly record is inserted as row 177
# (pushing the previously-existing rows 177-192 down to 178-193).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 177, shifting rows 177:192 down to 178:193.
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with the new weekly data point.
$ws.Cells.Item(177, 1).Value = 4
$ws.Cells.Item(177, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(177, 3).Value = "Los Lagos"
$ws.Cells.Item(177, 4).Value = 45265
$ws.Cells.Item(177, 5).Value = 10
$ws.Cells.Item(177, 6).Value = 100112026
$ws.Cells.Item(177, 7).Value = "Haba"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 160
$ws.Cells.Item(177, 11).Value = 16000
$ws.Cells.Item(177, 12).Value = 16000
$ws.Cells.Item(177, 13).Value = 16000
$ws.Cells.Item(177, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(177, 15).Value = "Región del Maule"
$ws.Cells.Item(177, 16).Value = 640
$ws.Cells.Item(177, 17).Value = 25
$ws.Cells.Item(177, 18).Value = "Hortaliza"
